# Weekly update: insert the newest price-report entries (Primera / Segunda)
# for "Agrícola del Norte S.A. de Arica - Caigua" as two new rows right
# below the header block (at rows 37-38), pushing all older rows down by
# two. The sheet's used range therefore grows from A1:R131 to A1:R133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at 37 and 38 (existing row 37 onward shifts
# down by two, automatically growing the sheet's dimension to row 133).
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(38).Insert()

# New row 37: Caigua, Primera, week of 2022-09-07 (serial 44811)
$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value = 44811
$ws.Cells.Item(37, 5).Value = 15
$ws.Cells.Item(37, 6).Value = 100112036
$ws.Cells.Item(37, 7).Value = "Caigua"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 130
$ws.Cells.Item(37, 11).Value = 7000
$ws.Cells.Item(37, 12).Value = 8000
$ws.Cells.Item(37, 13).Value = 7500
$ws.Cells.Item(37, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 375
$ws.Cells.Item(37, 17).Value = 20
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# New row 38: Caigua, Segunda, week of 2022-09-07 (serial 44811)
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(38, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value = 44811
$ws.Cells.Item(38, 5).Value = 15
$ws.Cells.Item(38, 6).Value = 100112036
$ws.Cells.Item(38, 7).Value = "Caigua"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Segunda"
$ws.Cells.Item(38, 10).Value = 160
$ws.Cells.Item(38, 11).Value = 6000
$ws.Cells.Item(38, 12).Value = 7000
$ws.Cells.Item(38, 13).Value = 6500
$ws.Cells.Item(38, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(38, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value = 325
$ws.Cells.Item(38, 17).Value = 20
$ws.Cells.Item(38, 18).Value = "Hortaliza"
